$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.9999999990259346
$ws.Range("C2").Value = 0.0000000001118007144738666
$ws.Range("D2").Value = 0.0000000001737751632210285
$ws.Range("E2").Value = 0.0000000006884875831429895
$ws.Range("F2").Value = 46066
